# Apply "Ajout des maquettes du jeu" edits to the Journal de travail workbook.

$wb = $excel.ActiveWorkbook

$sheetMock = $wb.Worksheets.Item("Feuil1")   # sheet1.xml - main "Journal de travail" table
$sheetStat = $wb.Worksheets.Item("Sheet1")   # sheet2.xml - statistics / chart sheet

# ---------------------------------------------------------------------------
# 1. Update / extend the "maquettes" (mockups) rows in the Journal table.
# ---------------------------------------------------------------------------

# Rows 30-34 are brand-new rows describing the various mockups realised.
# (Row 30's description is entered first so the new shared-string table keeps
# the same ordering as the authored workbook.)
$sheetMock.Range("B30").Value = 45415
$sheetMock.Range("C30").Value = 0.51041666666666663
$sheetMock.Range("D30").Value = 0.52777777777777779
$sheetMock.Range("F30").Value = "Analyse"
$sheetMock.Range("G30").Value = "Réalisation de la maquette du menu principale"

# Row 29 already existed (partial row) - fill in its End time and change its
# description to the new, more specific wording.
$sheetMock.Range("D29").Value = 0.51041666666666663
$sheetMock.Range("G29").Value = "Création de la maquette du menu priciaple"

$sheetMock.Range("B31").Value = 45415
$sheetMock.Range("C31").Value = 0.52777777777777779
$sheetMock.Range("D31").Value = 0.58333333333333337
$sheetMock.Range("F31").Value = "Analyse"
$sheetMock.Range("G31").Value = "Création de la maquette du menu du chronomètre"

$sheetMock.Range("B32").Value = 45415
$sheetMock.Range("C32").Value = 0.58333333333333337
$sheetMock.Range("D32").Value = 0.62847222222222221
$sheetMock.Range("F32").Value = "Analyse"
$sheetMock.Range("G32").Value = "Création de la maquette d'une partie"

$sheetMock.Range("B33").Value = 45415
$sheetMock.Range("C33").Value = 0.63888888888888895
$sheetMock.Range("D33").Value = 0.66666666666666663
$sheetMock.Range("F33").Value = "Analyse"
$sheetMock.Range("G33").Value = "Création de la maquette pour la fin d'une partie"

$sheetMock.Range("B34").Value = 45415
$sheetMock.Range("C34").Value = 0.66666666666666663
$sheetMock.Range("D34").Value = 0.6875
$sheetMock.Range("F34").Value = "Documentation"
$sheetMock.Range("G34").Value = "Écriture des maquette dans le dossier de projet"

# ---------------------------------------------------------------------------
# 2. Misc. formatting cell: give I11 on the stats sheet the Hyperlink style.
# ---------------------------------------------------------------------------
$sheetStat.Range("I11").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 3. Page setup tweaks.
# ---------------------------------------------------------------------------
$sheetMock.PageSetup.Zoom = 59

$sheetStat.PageSetup.Orientation = 2  # xlLandscape
$sheetStat.PageSetup.Zoom = 98
$sheetStat.PageSetup.FitToPagesTall = $false

# ---------------------------------------------------------------------------
# 4. View changes: active sheet / selections / zoom.
# ---------------------------------------------------------------------------
$sheetMock.Select()
$sheetMock.Range("A13").Select()
$sheetMock.Application.ActiveWindow.Zoom = 100

$sheetStat.Select()
$sheetStat.Range("N26").Select()
$sheetStat.Application.ActiveWindow.Zoom = 85

$wb.Save()
